$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B83: it should be a numeric value (3) instead of a text/inline string
$ws.Range("B83").Value = 3

# Add new row 84 with the new annotation entry
$ws.Range("A84").Value = "Ying Tang"

# B84 must remain a text value "4" (not auto-converted to a number)
$b84 = $ws.Range("B84")
$b84.NumberFormat = "@"
$b84.Value = "4"
$b84.ClearFormats()

$ws.Range("C84").Value = "Per your suggestions"
$ws.Range("D84").Value = "ACK"
$ws.Range("E84").Value = "EXP"
$ws.Range("F84").Value = "23ce80a1-f5c9-4d52-8c77-e985ea50fb2a"
$ws.Range("G84").Value = "H1uR4GZRZ_annotated.xlsx"
$ws.Range("H84").Value = "Per your suggestions, we have improved the draft by running additional experiments."
